$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'265.69"
$ws.Range("D3").Value = "'22.76"
$ws.Range("D4").Value = "'6.284"
$ws.Range("D5").Value = "'0.06156"
$ws.Range("D6").Value = "'3.586"
$ws.Range("D7").Value = "'6.704"
$ws.Range("D8").Value = "'1.344"
$ws.Range("D9").Value = "'0.8289"
$ws.Range("D10").Value = "'0.01357"
$ws.Range("D11").Value = "'0.1580"
$ws.Range("D12").Value = "'0.08238"
$ws.Range("D13").Value = "'0.03420"
$ws.Range("D14").Value = "'0.03139"
$ws.Range("D15").Value = "'0.09241"
$ws.Range("D16").Value = "'3.891"
$ws.Range("D18").Value = "'0.04831"
$ws.Range("D19").Value = "'0.006279"
$ws.Range("D20").Value = "'0.005264"
$ws.Range("B21").Value = "UpBots"
$ws.Range("C21").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D21").Value = "'0.007489"
$ws.Range("E21").Value = "20UpBotsUBXTBestin24h"
$ws.Range("B22").Value = "BitKan"
$ws.Range("C22").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D22").Value = "'0.001091"
$ws.Range("E22").Value = "21BitKanKAN"
$ws.Range("B23").Value = "NitroEx"
$ws.Range("C23").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D23").Value = "'0.0001500"
$ws.Range("E23").Value = "22NitroExNTX"
$ws.Range("B24").Value = "LEO"
$ws.Range("C24").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D24").Value = "'3.769"
$ws.Range("E24").Value = "23LEOLEO"
$ws.Range("B25").Value = "BTSEToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D25").Value = "'2.282"
$ws.Range("E25").Value = "24BTSETokenBTSE"
$ws.Range("B26").Value = "BitpandaEcosystemToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D26").Value = "'0.3379"
$ws.Range("E26").Value = "25BitpandaEcosystemTokenBEST"
$ws.Range("B27").Value = "ProBitToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D27").Value = "'0.1227"
$ws.Range("E27").Value = "26ProBitTokenPROB"
$ws.Range("D40").Value = "'0.04636"
$ws.Range("D41").Value = "'0.006954"
$ws.Range("D42").Value = "'0.1138"
$ws.Range("D44").Value = "'0.01075"
$ws.Range("D45").Value = "'0.00006159"
$ws.Range("D47").Value = "'0.7779"
$ws.Range("D48").Value = "'0.1973"
$ws.Range("E48").Value = "47BOLOBOLO"
